$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit shifts the dated price observations down by one row (rows 52-203),
# inserting a brand new observation at row 52 (D52/J52) and pushing what used to
# be the last row (203) down into a newly appended row 204.

$firstRow = 52
$lastRow = 203

# 1. Capture the "old" values for the columns that vary per row (D, J, K, M, P)
#    before we start overwriting anything.
$oldD = @{}
$oldJ = @{}
$oldK = @{}
$oldM = @{}
$oldP = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $oldD[$r] = $ws.Range("D$r").Value2()
    $oldJ[$r] = $ws.Range("J$r").Value2()
    $oldK[$r] = $ws.Range("K$r").Value2()
    $oldM[$r] = $ws.Range("M$r").Value2()
    $oldP[$r] = $ws.Range("P$r").Value2()
}

# 2. Shift rows 53..203 down by one: new row r gets old row (r-1)'s values.
for ($r = $lastRow; $r -ge ($firstRow + 1); $r--) {
    $src = $r - 1
    $ws.Range("D$r").Value = $oldD[$src]
    $ws.Range("D$r").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("J$r").Value = $oldJ[$src]
    $ws.Range("K$r").Value = $oldK[$src]
    $ws.Range("M$r").Value = $oldM[$src]
    $ws.Range("P$r").Value = $oldP[$src]
}

# 3. Row 52 becomes a brand new observation (same K/L/M/N/O/P/Q/R context as before,
#    only D and J change).
$ws.Range("D52").Value = 44560
$ws.Range("D52").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J52").Value = 3080

# 4. Append a new row 204 that carries what used to be row 203's data, reusing the
#    constant columns (A, B, C, E, F, G, H, I, L, N, O, Q, R) shared by every row.
$ws.Range("A204").Value = $ws.Range("A203").Value2()
$ws.Range("B204").Value = $ws.Range("B203").Value2()
$ws.Range("C204").Value = $ws.Range("C203").Value2()
$ws.Range("D204").Value = $oldD[$lastRow]
$ws.Range("D204").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E204").Value = $ws.Range("E203").Value2()
$ws.Range("F204").Value = $ws.Range("F203").Value2()
$ws.Range("G204").Value = $ws.Range("G203").Value2()
$ws.Range("H204").Value = $ws.Range("H203").Value2()
$ws.Range("I204").Value = $ws.Range("I203").Value2()
$ws.Range("J204").Value = $oldJ[$lastRow]
$ws.Range("K204").Value = $oldK[$lastRow]
$ws.Range("L204").Value = $ws.Range("L203").Value2()
$ws.Range("M204").Value = $oldM[$lastRow]
$ws.Range("N204").Value = $ws.Range("N203").Value2()
$ws.Range("O204").Value = $ws.Range("O203").Value2()
$ws.Range("P204").Value = $oldP[$lastRow]
$ws.Range("Q204").Value = $ws.Range("Q203").Value2()
$ws.Range("R204").Value = $ws.Range("R203").Value2()
